$d = $word.ActiveDocument

# 1. Extend the repository URL to point at the week8 homework subfolder.
$d.Content.Find.Execute("https://github.com/LeoSuzu/Data_Structure_and_Algorythms.git", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://github.com/LeoSuzu/Data_Structure_and_Algorythms/tree/main/Homeworks/week8", 2)

# 2. Shrink that paragraph's text (and paragraph mark) to 10pt (sz/szCs = 20 half-points).
$urlParagraph = $d.Paragraphs(3)
$urlParagraph.Range.Font.Size = 10
$urlParagraph.Range.Font.SizeBi = 10

# 3. Remove the now-unwanted trailing empty paragraph that followed the link
#    (delete from the end of the link paragraph's mark through the end of
#    the empty paragraph's mark, merging the two paragraph marks into one).
$trailingParagraph = $d.Paragraphs(4)
$deleteStart = $urlParagraph.Range.End - 1
$deleteEnd = $trailingParagraph.Range.End
$d.Range($deleteStart, $deleteEnd).Delete()
